$wb = $excel.ActiveWorkbook

# --- survey sheet: collapse the two begin group/begin table/.../end group
# blocks (rows 11-21, plus the blank row 10 separator) into a single
# "text | select | Table" row, now at row 10. ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("A10:A21").EntireRow.Delete()
$survey.Range("A10").Value = "text"
$survey.Range("B10").Value = "select"
$survey.Range("C10").Value = "Table"

# --- settings sheet: form_id changes from "Justtest-date" to "Justtest" ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("B2").Value = "Justtest"

# --- selection / active-tab bookkeeping: "settings" becomes the active
# sheet/tab, "survey" keeps its selection but is no longer the tab shown. ---
$survey.Range("B10").Select() | Out-Null
$settings.Activate() | Out-Null
$settings.Range("B2").Select() | Out-Null
